$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Force the cell to be written as literal text, even when the
    # value looks numeric (e.g. "244.94"), then strip the temporary
    # Text number-format back off so no stray style is left behind.
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Range('D2').Value = '30.398.10'
$ws.Range('E2').Value = '  -1.19%  '

$ws.Range('D3').Value = '1.872.34'
$ws.Range('E3').Value = '  -0.99%  '

$ws.Range('E4').Value = '  +0.03%  '

Set-TextValue 'D5' '244.94'
$ws.Range('E5').Value = '  -1.86%  '

Set-TextValue 'D6' '0.9999'
$ws.Range('E6').Value = '  +0.03%  '

Set-TextValue 'D7' '0.4716'
$ws.Range('E7').Value = '  -1.02%  '

Set-TextValue 'D8' '0.2876'
$ws.Range('E8').Value = '  -2.11%  '

Set-TextValue 'D9' '0.06493'
$ws.Range('E9').Value = '  -0.73%  '

$ws.Range('E10').Value = '  -1.06%  '

Set-TextValue 'D11' '100.14'
$ws.Range('E11').Value = '  +2.80%  '

Set-TextValue 'D12' '0.07792'
$ws.Range('E12').Value = '  +0.44%  '

$ws.Range('D13').Value = '1.872.87'
$ws.Range('E13').Value = '  -0.93%  '

Set-TextValue 'D14' '0.7300'
$ws.Range('E14').Value = '  -1.22%  '

Set-TextValue 'D15' '5.174'
$ws.Range('E15').Value = '  -1.46%  '

Set-TextValue 'D16' '286.16'
$ws.Range('E16').Value = '  +0.97%  '

$ws.Range('D17').Value = '30.381.56'
$ws.Range('E17').Value = '  -1.40%  '

Set-TextValue 'D18' '13.11'
$ws.Range('E18').Value = '  -0.72%  '

Set-TextValue 'D19' '0.9997'
$ws.Range('E19').Value = '  -0.03%  '

Set-TextValue 'D20' '0.000007497'
$ws.Range('E20').Value = '  -1.13%  '

$ws.Range('D21').Value = '2.115.31'
$ws.Range('E21').Value = '  -1.00%  '

Set-TextValue 'D22' '5.337'
$ws.Range('E22').Value = '  -0.10%  '

Set-TextValue 'D23' '0.9998'
$ws.Range('E23').Value = '  +0.00%  '

$ws.Range('E24').Value = '  +1.08%  '

Set-TextValue 'D25' '163.15'
$ws.Range('E25').Value = '  -0.78%  '

Set-TextValue 'D26' '9.042'
$ws.Range('E26').Value = '  -2.35%  '

Set-TextValue 'D27' '18.98'
$ws.Range('E27').Value = '  +0.11%  '

Set-TextValue 'D28' '1.898'
$ws.Range('E28').Value = '  -1.77%  '

Set-TextValue 'D29' '0.09678'
$ws.Range('E29').Value = '  -0.82%  '

Set-TextValue 'D30' '1.319'
$ws.Range('E30').Value = '  -1.98%  '

Set-TextValue 'D31' '1.488'
$ws.Range('E31').Value = '  -1.03%  '

Set-TextValue 'D32' '4.233'
$ws.Range('E32').Value = '  -1.84%  '

$ws.Range('E33').Value = '  -1.26%  '

$ws.Range('E34').Value = '  -1.74%  '

$ws.Range('E35').Value = '  -0.30%  '

Set-TextValue 'D36' '0.6889'
$ws.Range('E36').Value = '  -1.83%  '

Set-TextValue 'D37' '2.724'
$ws.Range('E37').Value = '  +0.09%  '

Set-TextValue 'D38' '0.01901'
$ws.Range('E38').Value = '  -0.98%  '

Set-TextValue 'D39' '2.847'
$ws.Range('E39').Value = '  +1.41%  '

Set-TextValue 'D40' '76.09'
$ws.Range('E40').Value = '  -0.14%  '

Set-TextValue 'D41' '6.283'
$ws.Range('E41').Value = '  -1.16%  '

Set-TextValue 'D42' '1.966'
$ws.Range('E42').Value = '  -3.25%  '

Set-TextValue 'D43' '0.4225'
$ws.Range('E43').Value = '  -1.04%  '

Set-TextValue 'D44' '0.9990'
$ws.Range('E44').Value = '  -0.07%  '

Set-TextValue 'D45' '0.8243'
$ws.Range('E45').Value = '  -1.87%  '

Set-TextValue 'D46' '101.20'
$ws.Range('E46').Value = '  -0.85%  '

Set-TextValue 'D47' '9.753'
$ws.Range('E47').Value = '  +3.29%  '

Set-TextValue 'D48' '7.021'
$ws.Range('E48').Value = '  -1.08%  '

Set-TextValue 'D49' '35.04'
$ws.Range('E49').Value = '  -2.22%  '

$ws.Range('E50').Value = '  -0.20%  '

Set-TextValue 'D51' '890.15'
$ws.Range('E51').Value = '  -3.99%  '
